# Fourth Commit with Pos and Negative scenario
#
# - SignIn!C2 (hyperlinked e-mail) and CreateAccount!F2 both held the shared
#   string "testjaga006@gmail.com" (index 28). Re-typing that same new value
#   into *both* cells lets the engine rewrite that shared-string slot in
#   place as "testjaga007@gmail.com" instead of appending a brand new entry.
# - SignIn!D2 gets a new password value "jaga@1234" appended as a new shared
#   string, while CreateAccount!N2/O2 keep pointing at the original
#   "jaga@12345" string.
# - The active sheet moves from CreateAccount to SignIn, with SignIn's
#   selection landing on D2 and CreateAccount's selection moving to F7.

$wb  = $excel.ActiveWorkbook
$ws1 = $wb.Worksheets.Item("SignIn")
$ws2 = $wb.Worksheets.Item("CreateAccount")

# --- Update the e-mail address used on both sheets (rewrites shared string 28 in place) ---
$ws1.Range("C2").Value = "testjaga007@gmail.com"
$ws2.Range("F2").Value = "testjaga007@gmail.com"

# --- Update the SignIn password cell to the new value (adds a new shared string) ---
$ws1.Range("D2").Value = "jaga@1234"

# --- Update CreateAccount's own selection before leaving it, so it is no longer the active tab ---
$ws2.Range("F7").Select()

# --- Finally select SignIn!D2, making SignIn the active sheet/tab ---
$ws1.Range("D2").Select()
